# Update test cases 12.29
# On the "invalid_settings_parameter" sheet, insert a new negative test case
# row right after the header row: trying to expand a pool with PDs of
# different types should fail with "Fail to create Pool".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalid_settings_parameter")

# Insert a new blank row above the current row 2 (shifts existing rows down,
# inheriting row 1's formatting the way Excel's own Insert does).
$ws.Rows.Item(2).Insert()

# Fill in the new row's three columns (send command / expect result / checkpoint).
# Order matters here only insofar as it controls the order new shared strings
# are appended in -- match author's original order (B, then A, then C).
$ws.Range("B2").Value = "different types of pd can not compose the one pool "
$ws.Range("A2").Value = "pool -a expand -i 0 -p 6,16"
$ws.Range("C2").Value = "Fail to create Pool"

# Move the selection to A9, matching the saved view state.
$ws.Activate()
$ws.Range("A9").Select()
